# Selenium test fixture update (selenium -> 4.8): refresh the region/code
# lookup values in column B of the active sheet to the new codes that came
# out of the regenerated test data set. Only the values themselves change;
# the labels in column A and the sheet layout stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "YD"
$ws.Range("B3").Value  = "TZ"
$ws.Range("B5").Value  = "O4"
$ws.Range("B6").Value  = "DV"
$ws.Range("B7").Value  = "GL1"
$ws.Range("B8").Value  = "XI"
$ws.Range("B9").Value  = "XK"
$ws.Range("B12").Value = "ST"
$ws.Range("B13").Value = "T1"
